$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value2 = "29.372.97"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value2 = "  -0.36%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value2 = "1.848.27"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value2 = "  -0.24%  "

# Row 4
$ws.Cells.Item(4, 5).Value2 = "  +0.06%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value2 = "240.69"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value2 = "  -0.11%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value2 = "0.6290"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value2 = "  -0.24%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value2 = "1.001"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value2 = "  +0.02%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value2 = "0.07587"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value2 = "  -1.24%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value2 = "0.2912"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value2 = "  -0.85%  "

# Row 10
$ws.Cells.Item(10, 5).Value2 = "  -0.93%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value2 = "0.07757"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value2 = "  +0.08%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value2 = "1.847.58"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value2 = "  -0.83%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value2 = "5.012"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value2 = "  -0.54%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value2 = "0.6781"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value2 = "  -0.42%  "

# Row 15
$ws.Cells.Item(15, 5).Value2 = "  -1.94%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value2 = "83.12"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value2 = "  -0.71%  "

# Row 17
$ws.Cells.Item(17, 2).Value2 = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(17, 3).Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value2 = "2.094.16"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value2 = "  -1.35%  "

# Row 18
$ws.Cells.Item(18, 2).Value2 = "Uniswap"
$ws.Cells.Item(18, 3).Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value2 = "6.111"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value2 = "  -1.19%  "

# Row 19
$ws.Cells.Item(19, 2).Value2 = "WrappedBTC"
$ws.Cells.Item(19, 3).Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value2 = "29.369.71"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value2 = "  -0.45%  "

# Row 20
$ws.Cells.Item(20, 2).Value2 = "BitcoinCash"
$ws.Cells.Item(20, 3).Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value2 = "229.37"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value2 = "  +0.14%  "

# Row 21
$ws.Cells.Item(21, 2).Value2 = "Avalanche"
$ws.Cells.Item(21, 3).Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value2 = "12.32"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value2 = "  -1.17%  "

# Row 22
$ws.Cells.Item(22, 2).Value2 = "Dai"
$ws.Cells.Item(22, 3).Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value2 = "1.001"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value2 = "  +0.05%  "

# Row 23
$ws.Cells.Item(23, 2).Value2 = "Chainlink"
$ws.Cells.Item(23, 3).Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value2 = "7.441"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value2 = "  -0.35%  "

# Row 24
$ws.Cells.Item(24, 2).Value2 = "BinanceUSD"
$ws.Cells.Item(24, 3).Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value2 = "1.001"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value2 = "  -0.02%  "

# Row 25
$ws.Cells.Item(25, 2).Value2 = "Monero"
$ws.Cells.Item(25, 3).Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value2 = "159.07"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value2 = "  +1.16%  "

# Row 26
$ws.Cells.Item(26, 2).Value2 = "Stellar"
$ws.Cells.Item(26, 3).Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value2 = "0.1396"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value2 = "  +0.96%  "

# Row 27
$ws.Cells.Item(27, 2).Value2 = "Cosmos"
$ws.Cells.Item(27, 3).Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value2 = "8.440"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value2 = "  +0.29%  "

# Row 28
$ws.Cells.Item(28, 2).Value2 = "EthereumClassic"
$ws.Cells.Item(28, 3).Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value2 = "17.69"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value2 = "  -0.11%  "

# Row 29
$ws.Cells.Item(29, 2).Value2 = "Toncoin"
$ws.Cells.Item(29, 3).Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value2 = "1.403"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value2 = "  +5.21%  "

# Row 30
$ws.Cells.Item(30, 2).Value2 = "PancakeSwap"
$ws.Cells.Item(30, 3).Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value2 = "1.470"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value2 = "  +0.13%  "

# Row 31
$ws.Cells.Item(31, 2).Value2 = "Hedera"
$ws.Cells.Item(31, 3).Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value2 = "0.05685"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value2 = "  -0.12%  "

# Row 32
$ws.Cells.Item(32, 2).Value2 = "Filecoin"
$ws.Cells.Item(32, 3).Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value2 = "4.110"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value2 = "  -0.66%  "

# Row 33
$ws.Cells.Item(33, 2).Value2 = "InternetComputer(DFINITY)"
$ws.Cells.Item(33, 3).Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value2 = "4.043"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value2 = "  -0.08%  "

# Row 34
$ws.Cells.Item(34, 2).Value2 = "ARBITRUM"
$ws.Cells.Item(34, 3).Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value2 = "1.154"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value2 = "  -1.14%  "

# Row 35
$ws.Cells.Item(35, 2).Value2 = "LidoDAOToken"
$ws.Cells.Item(35, 3).Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value2 = "1.822"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value2 = "  -1.86%  "

# Row 36
$ws.Cells.Item(36, 2).Value2 = "ImmutableX"
$ws.Cells.Item(36, 3).Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value2 = "0.6970"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value2 = "  -1.64%  "

# Row 37
$ws.Cells.Item(37, 2).Value2 = "HuobiToken"
$ws.Cells.Item(37, 3).Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value2 = "2.584"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value2 = "  -0.16%  "

# Row 38
$ws.Cells.Item(38, 2).Value2 = "VeChain"
$ws.Cells.Item(38, 3).Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value2 = "0.01832"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value2 = "  +2.06%  "

# Row 39
$ws.Cells.Item(39, 2).Value2 = "Maker"
$ws.Cells.Item(39, 3).Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value2 = "1.241.96"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value2 = "  +1.73%  "

# Row 40
$ws.Cells.Item(40, 2).Value2 = "MXToken"
$ws.Cells.Item(40, 3).Value2 = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value2 = "2.719"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value2 = "  -2.39%  "

# Row 41
$ws.Cells.Item(41, 2).Value2 = "FraxShare"
$ws.Cells.Item(41, 3).Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value2 = "6.405"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value2 = "  -2.39%  "

# Row 42
$ws.Cells.Item(42, 2).Value2 = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value2 = "0.9019"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value2 = "  -0.70%  "

# Row 43
$ws.Cells.Item(43, 2).Value2 = "PaxDollar"
$ws.Cells.Item(43, 3).Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value2 = "1.000"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value2 = "  -0.10%  "

# Row 44
$ws.Cells.Item(44, 2).Value2 = "RocketPoolETH"
$ws.Cells.Item(44, 3).Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value2 = "2.003.71"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value2 = "  -1.33%  "

# Row 45
$ws.Cells.Item(45, 2).Value2 = "Quant"
$ws.Cells.Item(45, 3).Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value2 = "101.50"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value2 = "  -0.27%  "

# Row 46
$ws.Cells.Item(46, 2).Value2 = "Aave"
$ws.Cells.Item(46, 3).Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value2 = "65.51"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value2 = "  -1.62%  "

# Row 47
$ws.Cells.Item(47, 2).Value2 = "Aptos"
$ws.Cells.Item(47, 3).Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value2 = "7.127"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value2 = "  -0.05%  "

# Row 48
$ws.Cells.Item(48, 2).Value2 = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value2 = "0.00000000118"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value2 = "  -2.35%  "

# Row 49
$ws.Cells.Item(49, 2).Value2 = "TheSandbox"
$ws.Cells.Item(49, 3).Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value2 = "0.3996"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value2 = "  -0.73%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value2 = "0.1151"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value2 = "  +0.50%  "

# Row 51
$ws.Cells.Item(51, 2).Value2 = "EnergySwap"
$ws.Cells.Item(51, 3).Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value2 = "8.986"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value2 = "  -0.39%  "
